# Sample.xlsx "Add files via upload" re-edit
#
# 1) The two wrapped quarter headers ("YTD\nQ2'16" / "YTD\nQ3'16") lose
#    their embedded line break and become single-line labels
#    ("YTD Q2'16" / "YTD Q3'16") on both worksheets. Write R1 before L1
#    so the shared-string table picks up "YTD Q3'16" ahead of
#    "YTD Q2'16", matching the upload's string order.
# 2) The active tab moves from Sheet1 to Sheet2, and the remembered
#    selection on each sheet changes (Sheet1 -> X1, Sheet2 -> P1).

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("Sheet1", "Sheet2")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("R1").Value = "YTD Q3'16"
    $ws.Range("L1").Value = "YTD Q2'16"
}

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("X1").Select()

$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Activate()
$ws2.Range("P1").Select()
